$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2830
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F5").Value = 1576
$ws1.Range("F6").Value = 1157
$ws1.Range("F7").Value = 328
$ws1.Range("F12").Value = 9542
$ws1.Range("F13").Value = 409
$ws1.Range("F14").Value = 2513
$ws1.Range("F15").Value = 12
$ws1.Range("F16").Value = 269
$ws1.Range("F17").Value = 185
$ws1.Range("F18").Value = 474
$ws1.Range("F19").Value = 685
$ws1.Range("F20").Value = 685
$ws1.Range("F21").Value = 1197
$ws1.Range("F22").Value = 1005
$ws1.Range("F23").Value = 2968
$ws1.Range("F24").Value = 2245
$ws1.Range("F25").Value = 1935
$ws1.Range("F29").Value = 1558
$ws1.Range("F31").Value = 17
$ws1.Range("F32").Value = 177
$ws1.Range("F35").Value = 343
$ws1.Range("F36").Value = 69
$ws1.Range("F38").Value = 511
$ws1.Range("F40").Value = 123
$ws1.Range("F41").Value = 1538
$ws1.Range("F42").Value = 133
$ws1.Range("F43").Value = 1485
$ws1.Range("F44").Value = 30
$ws1.Range("F45").Value = 342
$ws1.Range("F46").Value = 24
$ws1.Range("F47").Value = 365
$ws1.Range("F48").Value = 745
$ws1.Range("F50").Value = 314

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2830
$ws4.Range("F4").Value = 1576
$ws4.Range("F6").Value = 1157
$ws4.Range("F8").Value = 9543
$ws4.Range("F9").Value = 409
$ws4.Range("F10").Value = 2513
$ws4.Range("F12").Value = 12
$ws4.Range("F14").Value = 269
$ws4.Range("F15").Value = 185
$ws4.Range("F16").Value = 685
$ws4.Range("F17").Value = 685
$ws4.Range("F18").Value = 1197
$ws4.Range("F19").Value = 1005
$ws4.Range("F20").Value = 2968
$ws4.Range("F21").Value = 2245
$ws4.Range("F22").Value = 1935
$ws4.Range("F24").Value = 1558
$ws4.Range("F26").Value = 17
$ws4.Range("F27").Value = 177
$ws4.Range("F30").Value = 343
$ws4.Range("F31").Value = 69
$ws4.Range("F33").Value = 511
$ws4.Range("F38").Value = 123
$ws4.Range("F39").Value = 1538
$ws4.Range("F41").Value = 133
$ws4.Range("F42").Value = 1485
$ws4.Range("F43").Value = 30
$ws4.Range("F45").Value = 342
$ws4.Range("F46").Value = 24
$ws4.Range("F47").Value = 365
$ws4.Range("F48").Value = 745
$ws4.Range("F49").Value = 314
